$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 231.5
$ws.Range("I2").Value = 72.75
$ws.Range("J2").Value = 549
$ws.Range("K2").Value = 72.75
$ws.Range("L2").Value = 549
$ws.Range("M2").Value = 40.25
$ws.Range("N2").Value = -775
$ws.Range("H40").Value = 1490
$ws.Range("I40").Value = 1416.2069
$ws.Range("J40").Value = 2025
$ws.Range("K40").Value = 1416.2069
$ws.Range("L40").Value = 2025
$ws.Range("M40").Value = -1241.2069
$ws.Range("N40").Value = -2375
$ws.Range("H137").Value = 1336.129
$ws.Range("I137").Value = 1176.4584
$ws.Range("J137").Value = 1883.5714
$ws.Range("K137").Value = 3529.3752
$ws.Range("L137").Value = 5650.7142
$ws.Range("M137").Value = -979.3751999999999
$ws.Range("N137").Value = -10750.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9720.416999999999
$ws.Range("I45").Value = 11364.9
$ws.Range("K45").Value = 11364.9
$ws.Range("M45").Value = -10987.9
$ws.Range("H61").Value = 10192.77
$ws.Range("I61").Value = 10875.5
$ws.Range("K61").Value = 10875.5
$ws.Range("M61").Value = -10663.5
$ws.Range("H74").Value = 2898.4614
$ws.Range("I74").Value = 3342
$ws.Range("J74").Value = 2621.25
$ws.Range("K74").Value = 3342
$ws.Range("L74").Value = 2621.25
$ws.Range("M74").Value = -2468
$ws.Range("N74").Value = -4369.25
$ws.Range("H77").Value = 2898.4614
$ws.Range("I77").Value = 3342
$ws.Range("J77").Value = 2621.25
$ws.Range("K77").Value = 16710
$ws.Range("L77").Value = 13106.25
$ws.Range("M77").Value = -12342
$ws.Range("N77").Value = -21842.25
$ws.Range("H132").Value = 2729.1282
$ws.Range("I132").Value = 1261.6316
$ws.Range("J132").Value = 4123.25
$ws.Range("K132").Value = 3784.8948
$ws.Range("L132").Value = 12369.75
$ws.Range("M132").Value = -1254.8948
$ws.Range("N132").Value = -17429.75
$ws.Range("H136").Value = 10192.77
$ws.Range("I136").Value = 10875.5
$ws.Range("K136").Value = 32626.5
$ws.Range("M136").Value = -30076.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1356.8572
$ws.Range("J107").Value = 999
$ws.Range("L107").Value = 999
$ws.Range("N107").Value = -4839
$ws.Range("H134").Value = 5075.769
$ws.Range("I134").Value = 7138.8096
$ws.Range("K134").Value = 21416.4288
$ws.Range("M134").Value = -18881.4288

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 143.125
$ws.Range("I7").Value = 152.6
$ws.Range("J7").Value = 127.333336
$ws.Range("K7").Value = 152.6
$ws.Range("L7").Value = 127.333336
$ws.Range("M7").Value = -39.59999999999999
$ws.Range("N7").Value = -353.333336
$ws.Range("H22").Value = 513.2353000000001
$ws.Range("I22").Value = 448.07693
$ws.Range("K22").Value = 448.07693
$ws.Range("M22").Value = -98.07693
$ws.Range("H31").Value = 6078.5835
$ws.Range("I31").Value = 1896.091
$ws.Range("J31").Value = 12651.071
$ws.Range("K31").Value = 1896.091
$ws.Range("L31").Value = 12651.071
$ws.Range("M31").Value = -1601.091
$ws.Range("N31").Value = -13241.071
$ws.Range("H34").Value = 6078.5835
$ws.Range("I34").Value = 1896.091
$ws.Range("J34").Value = 12651.071
$ws.Range("K34").Value = 1896.091
$ws.Range("L34").Value = 12651.071
$ws.Range("M34").Value = -1694.091
$ws.Range("N34").Value = -13055.071
$ws.Range("H58").Value = 1382.9375
$ws.Range("I58").Value = 1082.25
$ws.Range("J58").Value = 1884.0834
$ws.Range("K58").Value = 1082.25
$ws.Range("L58").Value = 1884.0834
$ws.Range("M58").Value = -879.25
$ws.Range("N58").Value = -2290.0834
$ws.Range("H122").Value = 1998
$ws.Range("I122").Value = 1997.75
$ws.Range("J122").Value = 1998.2
$ws.Range("K122").Value = 5993.25
$ws.Range("L122").Value = 5994.6
$ws.Range("M122").Value = -3543.25
$ws.Range("N122").Value = -10894.6
$ws.Range("H132").Value = 3350.9565
$ws.Range("I132").Value = 3360.1177
$ws.Range("K132").Value = 10080.3531
$ws.Range("M132").Value = -7550.3531
$ws.Range("H134").Value = 3351.64
$ws.Range("I134").Value = 4016.6316
$ws.Range("J134").Value = 1245.8334
$ws.Range("K134").Value = 12049.8948
$ws.Range("L134").Value = 3737.5002
$ws.Range("M134").Value = -9514.8948
$ws.Range("N134").Value = -8807.5002
$ws.Range("H136").Value = 1382.9375
$ws.Range("I136").Value = 1082.25
$ws.Range("J136").Value = 1884.0834
$ws.Range("K136").Value = 3246.75
$ws.Range("L136").Value = 5652.2502
$ws.Range("M136").Value = -696.75
$ws.Range("N136").Value = -10752.2502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 500
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 52.4
$ws.Range("I2").Value = 14
$ws.Range("J2").Value = 110
$ws.Range("K2").Value = 14
$ws.Range("L2").Value = 110
$ws.Range("M2").Value = 99
$ws.Range("N2").Value = -336
$ws.Range("H132").Value = 5432.8125
$ws.Range("I132").Value = 6309.9
$ws.Range("K132").Value = 18929.7
$ws.Range("M132").Value = -16399.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1950800.1
$ws.Range("I22").Value = 5291466
$ws.Range("J22").Value = 2078.5278
$ws.Range("K22").Value = 5291466
$ws.Range("L22").Value = 2078.5278
$ws.Range("M22").Value = -5291171
$ws.Range("N22").Value = -2668.5278
$ws.Range("H27").Value = 1950800.1
$ws.Range("I27").Value = 5291466
$ws.Range("J27").Value = 2078.5278
$ws.Range("K27").Value = 5291466
$ws.Range("L27").Value = 2078.5278
$ws.Range("M27").Value = -5291359
$ws.Range("N27").Value = -2292.5278
$ws.Range("H40").Value = 166668780
$ws.Range("I40").Value = 200001950
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 200001950
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -200001814
$ws.Range("N40").Value = -3272
$ws.Range("H68").Value = 43480440
$ws.Range("I68").Value = 1698.3334
$ws.Range("J68").Value = 90911790
$ws.Range("K68").Value = 1698.3334
$ws.Range("L68").Value = 90911790
$ws.Range("M68").Value = -949.3334
$ws.Range("N68").Value = -90913288
$ws.Range("H71").Value = 43480440
$ws.Range("I71").Value = 1698.3334
$ws.Range("J71").Value = 90911790
$ws.Range("K71").Value = 8491.666999999999
$ws.Range("L71").Value = 454558950
$ws.Range("M71").Value = -4747.666999999999
$ws.Range("N71").Value = -454566438
$ws.Range("H132").Value = 14790679
$ws.Range("I132").Value = 24130286
$ws.Range("J132").Value = 2966.0833
$ws.Range("K132").Value = 72390858
$ws.Range("L132").Value = 8898.249899999999
$ws.Range("M132").Value = -72388328
$ws.Range("N132").Value = -13958.2499
$ws.Range("H136").Value = 6780.846
$ws.Range("I136").Value = 6644.5557
$ws.Range("K136").Value = 19933.6671
$ws.Range("M136").Value = -17383.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 15000
$ws.Range("J47").Value = 15000
$ws.Range("L47").Value = 15000
$ws.Range("N47").Value = -16144
$ws.Range("H122").Value = 3998.8
$ws.Range("I122").Value = 3331.8333
$ws.Range("K122").Value = 9995.499899999999
$ws.Range("M122").Value = -7545.499899999999
$ws.Range("H132").Value = 1982.091
$ws.Range("I132").Value = 1266.5
$ws.Range("J132").Value = 2840.8
$ws.Range("K132").Value = 3799.5
$ws.Range("L132").Value = 8522.400000000001
$ws.Range("M132").Value = -1269.5
$ws.Range("N132").Value = -13582.4
$ws.Range("H136").Value = 2215.9524
$ws.Range("I136").Value = 2617.8
$ws.Range("J136").Value = 1625
$ws.Range("K136").Value = 7853.400000000001
$ws.Range("L136").Value = 4875
$ws.Range("M136").Value = -5303.400000000001
$ws.Range("N136").Value = -9975
